$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L54").ClearContents()
$ws.Range("M54").ClearContents()
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
$ws.Range("H9").Value = 153.57143
$ws.Range("I9").Value = 145
$ws.Range("J9").Value = 157
$ws.Range("K9").Value = 145
$ws.Range("L9").Value = 157
$ws.Range("M9").Value = 24
$ws.Range("N9").Value = -495
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 2000
$ws.Range("J17").Value = 6000
$ws.Range("K17").Value = 6000
$ws.Range("L17").Value = 18000
$ws.Range("M17").Value = -5832
$ws.Range("N17").Value = -18336
$ws.Range("H19").Value = 144
$ws.Range("J19").Value = 14
$ws.Range("L19").Value = 14
$ws.Range("N19").Value = -364
$ws.Range("H51").Value = 4997
$ws.Range("J51").Value = 4996.3335
$ws.Range("L51").Value = 4996.3335
$ws.Range("N51").Value = -5964.3335
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("N54").Value = 0
$ws.Range("H62").Value = 4941.5386
$ws.Range("I62").Value = 4941.5386
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4941.5386
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = -4317.5386
$ws.Range("H65").Value = 4941.5386
$ws.Range("I65").Value = 4941.5386
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 24707.693
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = -21587.693
$ws.Range("H69").Value = 9999
$ws.Range("J69").Value = 9999
$ws.Range("L69").Value = 29997
$ws.Range("N69").Value = -31745
$ws.Range("H72").Value = 9999
$ws.Range("J72").Value = 9999
$ws.Range("L72").Value = 89991
$ws.Range("N72").Value = -98727
$ws.Range("H98").Value = 1831.6471
$ws.Range("I98").Value = 591.3
$ws.Range("K98").Value = 591.3
$ws.Range("M98").Value = 906.7
$ws.Range("H122").Value = 1831.6471
$ws.Range("I122").Value = 591.3
$ws.Range("K122").Value = 1773.9
$ws.Range("M122").Value = 676.1000000000001
$ws.Range("H132").Value = 30306438
$ws.Range("I132").Value = 35717604
$ws.Range("K132").Value = 107152812
$ws.Range("M132").Value = -107150282
$ws.Range("H137").Value = 1894.7858
$ws.Range("I137").Value = 1543.9166
$ws.Range("J137").Value = 4000
$ws.Range("K137").Value = 4631.7498
$ws.Range("L137").Value = 12000
$ws.Range("M137").Value = -2081.7498
$ws.Range("N137").Value = -17100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M102").ClearContents()
$ws.Range("H2").Value = 988.3333
$ws.Range("I2").Value = 988.3333
$ws.Range("K2").Value = 988.3333
$ws.Range("M2").Value = -875.3333
$ws.Range("H24").Value = 50177.5
$ws.Range("J24").Value = 50177.5
$ws.Range("L24").Value = 50177.5
$ws.Range("N24").Value = -50925.5
$ws.Range("H32").Value = 9166.462
$ws.Range("I32").Value = 9333.360000000001
$ws.Range("J32").Value = 4994
$ws.Range("K32").Value = 9333.360000000001
$ws.Range("L32").Value = 4994
$ws.Range("M32").Value = -9046.360000000001
$ws.Range("N32").Value = -5568
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
$ws.Range("H74").Value = 2413.4285
$ws.Range("I74").Value = 2179.8
$ws.Range("J74").Value = 2997.5
$ws.Range("K74").Value = 2179.8
$ws.Range("L74").Value = 2997.5
$ws.Range("M74").Value = -1305.8
$ws.Range("N74").Value = -4745.5
$ws.Range("H77").Value = 2413.4285
$ws.Range("I77").Value = 2179.8
$ws.Range("J77").Value = 2997.5
$ws.Range("K77").Value = 10899
$ws.Range("L77").Value = 14987.5
$ws.Range("M77").Value = -6531
$ws.Range("N77").Value = -23723.5
$ws.Range("H92").Value = 34498
$ws.Range("J92").Value = 34498
$ws.Range("L92").Value = 34498
$ws.Range("N92").Value = -39490
$ws.Range("H100").Value = 50177.5
$ws.Range("J100").Value = 50177.5
$ws.Range("L100").Value = 50177.5
$ws.Range("N100").Value = -52341.5
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("H116").Value = 988.3333
$ws.Range("I116").Value = 988.3333
$ws.Range("K116").Value = 988.3333
$ws.Range("M116").Value = 1305.6667
$ws.Range("H122").Value = 3767
$ws.Range("I122").Value = 3836.7273
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 11510.1819
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -9060.1819
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M105").ClearContents()
$ws.Range("H3").Value = 988.3333
$ws.Range("I3").Value = 988.3333
$ws.Range("K3").Value = 988.3333
$ws.Range("M3").Value = -874.3333
$ws.Range("H20").Value = 3249.1428
$ws.Range("I20").Value = 1949.75
$ws.Range("J20").Value = 4981.6665
$ws.Range("K20").Value = 1949.75
$ws.Range("L20").Value = 4981.6665
$ws.Range("M20").Value = -1702.75
$ws.Range("N20").Value = -5475.6665
$ws.Range("H86").Value = 2661.875
$ws.Range("I86").Value = 2661.875
$ws.Range("K86").Value = 2661.875
$ws.Range("M86").Value = -1538.875
$ws.Range("H89").Value = 2661.875
$ws.Range("I89").Value = 2661.875
$ws.Range("K89").Value = 13309.375
$ws.Range("M89").Value = -7693.375
$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("H107").Value = 1128.091
$ws.Range("I107").Value = 670
$ws.Range("K107").Value = 670
$ws.Range("M107").Value = 1250
$ws.Range("H130").Value = 64666.332
$ws.Range("J130").Value = 64666.332
$ws.Range("L130").Value = 64666.332
$ws.Range("N130").Value = -74706.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L140").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("H16").Value = 8139
$ws.Range("I16").Value = 3565
$ws.Range("K16").Value = 3565
$ws.Range("M16").Value = -3278
$ws.Range("H17").Value = 776.5
$ws.Range("I17").Value = 106
$ws.Range("K17").Value = 106
$ws.Range("M17").Value = 68
$ws.Range("H22").Value = 648.625
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150
$ws.Range("H26").Value = 3020.2
$ws.Range("I26").Value = 1101
$ws.Range("J26").Value = 3500
$ws.Range("K26").Value = 1101
$ws.Range("L26").Value = 3500
$ws.Range("M26").Value = -814
$ws.Range("N26").Value = -4074
$ws.Range("H28").Value = 22699.5
$ws.Range("J28").Value = 22699.5
$ws.Range("L28").Value = 22699.5
$ws.Range("N28").Value = -23189.5
$ws.Range("H31").Value = 4273.3076
$ws.Range("I31").Value = 4116.3335
$ws.Range("J31").Value = 4407.857
$ws.Range("K31").Value = 4116.3335
$ws.Range("L31").Value = 4407.857
$ws.Range("M31").Value = -3821.3335
$ws.Range("N31").Value = -4997.857
$ws.Range("H34").Value = 4273.3076
$ws.Range("I34").Value = 4116.3335
$ws.Range("J34").Value = 4407.857
$ws.Range("K34").Value = 4116.3335
$ws.Range("L34").Value = 4407.857
$ws.Range("M34").Value = -3914.3335
$ws.Range("N34").Value = -4811.857
$ws.Range("H58").Value = 3636.4
$ws.Range("I58").Value = 2394
$ws.Range("K58").Value = 2394
$ws.Range("M58").Value = -2191
$ws.Range("H62").Value = 6194.125
$ws.Range("I62").Value = 6194.125
$ws.Range("K62").Value = 6194.125
$ws.Range("M62").Value = -5570.125
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
$ws.Range("H65").Value = 6194.125
$ws.Range("I65").Value = 6194.125
$ws.Range("K65").Value = 30970.625
$ws.Range("M65").Value = -27850.625
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
$ws.Range("H82").Value = 22500
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 30000
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 30000
$ws.Range("M82").Value = -14639
$ws.Range("N82").Value = -30722
$ws.Range("H85").Value = 22500
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 30000
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 30000
$ws.Range("M85").Value = -13752
$ws.Range("N85").Value = -32496
$ws.Range("H99").Value = 2775
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 3433.3333
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 3433.3333
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -6429.3333
$ws.Range("H109").Value = 58659.332
$ws.Range("J109").Value = 58659.332
$ws.Range("L109").Value = 58659.332
$ws.Range("N109").Value = -60739.332
$ws.Range("H113").Value = 8139
$ws.Range("I113").Value = 3565
$ws.Range("K113").Value = 3565
$ws.Range("M113").Value = -1395
$ws.Range("H126").Value = 2775
$ws.Range("I126").Value = 800
$ws.Range("J126").Value = 3433.3333
$ws.Range("K126").Value = 2400
$ws.Range("L126").Value = 10299.9999
$ws.Range("M126").Value = 70
$ws.Range("N126").Value = -15239.9999
$ws.Range("H134").Value = 3927.5833
$ws.Range("I134").Value = 2813.5
$ws.Range("J134").Value = 9498
$ws.Range("K134").Value = 8440.5
$ws.Range("L134").Value = 28494
$ws.Range("M134").Value = -5905.5
$ws.Range("N134").Value = -33564
$ws.Range("H136").Value = 3636.4
$ws.Range("I136").Value = 2394
$ws.Range("K136").Value = 7182
$ws.Range("M136").Value = -4632
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").Value = 0
$ws.Range("H141").Value = 100000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = 100000
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L16").ClearContents()
$ws.Range("M139").ClearContents()
$ws.Range("H4").Value = 1692.25
$ws.Range("I4").Value = 1651.1818
$ws.Range("J4").Value = 1727
$ws.Range("K4").Value = 4953.5454
$ws.Range("L4").Value = 5181
$ws.Range("M4").Value = -4841.5454
$ws.Range("N4").Value = -5405
$ws.Range("H16").Value = 1000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 0
$ws.Range("M16").Value = 3000
$ws.Range("N16").Value = -3346
$ws.Range("H29").Value = 714
$ws.Range("I29").Value = 280
$ws.Range("J29").Value = 887.6
$ws.Range("K29").Value = 840
$ws.Range("L29").Value = 2662.8
$ws.Range("M29").Value = -563
$ws.Range("N29").Value = -3216.8
$ws.Range("H34").Value = 3120.75
$ws.Range("I34").Value = 466.66666
$ws.Range("J34").Value = 4713.2
$ws.Range("K34").Value = 1399.99998
$ws.Range("L34").Value = 14139.6
$ws.Range("M34").Value = -1315.99998
$ws.Range("N34").Value = -14307.6
$ws.Range("H92").Value = 500
$ws.Range("J92").Value = 500
$ws.Range("L92").Value = 1500
$ws.Range("N92").Value = -3996
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L133").ClearContents()
$ws.Range("H97").Value = 1103.4286
$ws.Range("I97").Value = 1120.75
$ws.Range("K97").Value = 1120.75
$ws.Range("M97").Value = -624.75
$ws.Range("H102").Value = 3208.6
$ws.Range("I102").Value = 3208.6
$ws.Range("K102").Value = 3208.6
$ws.Range("M102").Value = -1586.6
$ws.Range("H113").Value = 1616.3334
$ws.Range("I113").Value = 1616.3334
$ws.Range("K113").Value = 1616.3334
$ws.Range("M113").Value = 553.6666
$ws.Range("H132").Value = 989
$ws.Range("I132").Value = 989
$ws.Range("K132").Value = 2967
$ws.Range("M132").Value = -437
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 2982.6667
$ws.Range("J11").Value = 2982.6667
$ws.Range("L11").Value = 2982.6667
$ws.Range("N11").Value = -3262.6667
$ws.Range("H40").Value = 6000.2
$ws.Range("I40").Value = 5000.25
$ws.Range("K40").Value = 5000.25
$ws.Range("M40").Value = -4864.25
$ws.Range("H61").Value = 1000
$ws.Range("I61").Value = 1000
$ws.Range("K61").Value = 1000
$ws.Range("M61").Value = -798
$ws.Range("H82").Value = 1509.3077
$ws.Range("I82").Value = 1509.3077
$ws.Range("K82").Value = 1509.3077
$ws.Range("M82").Value = -1148.3077
$ws.Range("H85").Value = 1509.3077
$ws.Range("I85").Value = 1509.3077
$ws.Range("K85").Value = 1509.3077
$ws.Range("M85").Value = -261.3077000000001
$ws.Range("H113").Value = 1000
$ws.Range("I113").Value = 1000
$ws.Range("K113").Value = 1000
$ws.Range("M113").Value = 1170
$ws.Range("H136").Value = 3184.375
$ws.Range("I136").Value = 3184.375
$ws.Range("K136").Value = 9553.125
$ws.Range("M136").Value = -7003.125
$ws.Range("H140").Value = 125900
$ws.Range("J140").Value = 125900
$ws.Range("L140").Value = 125900
$ws.Range("N140").Value = -136260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1050
$ws.Range("H65").Value = 1050
$ws.Range("H81").Value = 2995
$ws.Range("I81").Value = 2995
$ws.Range("K81").Value = 5990
$ws.Range("M81").Value = -4929
$ws.Range("H82").Value = 7999
$ws.Range("J82").Value = 7999
$ws.Range("L82").Value = 7999
$ws.Range("N82").Value = -8765
$ws.Range("H84").Value = 2995
$ws.Range("I84").Value = 2995
$ws.Range("K84").Value = 29950
$ws.Range("M84").Value = -24646
$ws.Range("H85").Value = 7999
$ws.Range("J85").Value = 7999
$ws.Range("L85").Value = 7999
$ws.Range("N85").Value = -10651
$ws.Range("H107").Value = 349.2857
$ws.Range("I107").Value = 357.5
$ws.Range("K107").Value = 1072.5
$ws.Range("M107").Value = 847.5
$ws.Range("H122").Value = 1137.0667
$ws.Range("I122").Value = 1147.5714
$ws.Range("J122").Value = 990
$ws.Range("K122").Value = 3442.7142
$ws.Range("L122").Value = 2970
$ws.Range("M122").Value = -992.7142000000003
$ws.Range("N122").Value = -7870
